$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 2: only CV stat columns change ----
$ws.Range("H2").Value = 0.9706992789427162
$ws.Range("I2").Value = 0.007943953096139331
$ws.Range("J2").Value = 0.5185769661710837
$ws.Range("K2").Value = 0.1538086624801142

# ---- Row 3 ----
$A3 = @'
Pipeline(steps=[('scaler', RobustScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f91c9d63760>),
                ('model',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',
                                                                     max_depth=6,
                                                                     max_features='sqrt',
                                                                     min_samples_split=4,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])
'@
$ws.Range("A3").Value = $A3

$ws.Range("B3").Value = 0.699871794871795

$C3 = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91243515b0>, 'scaler': RobustScaler(), 'model__n_estimators': 10, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 1, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 6, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': 'balanced'}"
$ws.Range("C3").Value = $C3

$ws.Range("D3").Value = 0.75

$ws.Range("F3").Value = "[0 1 1 1 1 0 1 1 1 1 1 0]"

$ws.Range("H3").Value = 0.9724702579371853
$ws.Range("I3").Value = 0.007081119230557222
$ws.Range("J3").Value = 0.6025795465971937
$ws.Range("K3").Value = 0.1234478071525753

# ---- Row 4 ----
$A4 = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f91c7468e80>),
                ('model',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(max_depth=6,
                                                                     min_samples_leaf=5,
                                                                     min_samples_split=6,
                                                                     random_state=42),
                                    random_state=42))])
'@
$ws.Range("A4").Value = $A4

$ws.Range("B4").Value = 0.7034523809523809

$C4 = "{'selector': <__main__.NamedFeatureSelector object at 0x7f9124366970>, 'scaler': MinMaxScaler(), 'model__n_estimators': 50, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': None, 'model__estimator__max_depth': 6, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}"
$ws.Range("C4").Value = $C4

$ws.Range("D4").Value = 0.6153846153846154

$ws.Range("F4").Value = "[0 1 1 0 0 1 0 1 1 0 0 0]"

$ws.Range("H4").Value = 0.9739405914781303
$ws.Range("I4").Value = 0.006152057408022343
$ws.Range("J4").Value = 0.6035929364752894
$ws.Range("K4").Value = 0.1268847575029185

# ---- Row 5: only CV stat columns change ----
$ws.Range("H5").Value = 0.9753857809243834
$ws.Range("I5").Value = 0.005251340185646229
$ws.Range("J5").Value = 0.6038589580354286
$ws.Range("K5").Value = 0.1390832661757343
